$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated cell values. Price-column values that look numeric are
# prefixed with a literal apostrophe so Excel keeps storing them as text
# (matching the original inline-string "Price" column formatting) instead
# of silently converting them into numbers.
$ws.Range("D2").Value = '89.548.27'
$ws.Range("E2").Value = '  -1.60%  '
$ws.Range("D3").Value = '3.136.37'
$ws.Range("E3").Value = '  -2.49%  '
$ws.Range("D4").Value = '''1.00'
$ws.Range("E4").Value = '  -0.41%  '
$ws.Range("D5").Value = '''215.18'
$ws.Range("E5").Value = '  -0.31%  '
$ws.Range("D6").Value = '''638.21'
$ws.Range("E6").Value = '  +2.87%  '
$ws.Range("D7").Value = '''0.397'
$ws.Range("E7").Value = '  +2.06%  '
$ws.Range("D8").Value = '''0.771'
$ws.Range("E8").Value = '  +10.55%  '
$ws.Range("D9").Value = '''1.00'
$ws.Range("E9").Value = '  -0.08%  '
$ws.Range("D10").Value = '3.132.76'
$ws.Range("E10").Value = '  -2.03%  '
$ws.Range("E11").Value = '  -3.16%  '
$ws.Range("E12").Value = '  +0.44%  '
$ws.Range("E13").Value = '  -2.40%  '
$ws.Range("E14").Value = '  +1.73%  '
$ws.Range("D15").Value = '89.311.61'
$ws.Range("E15").Value = '  -1.87%  '
$ws.Range("D16").Value = '3.706.83'
$ws.Range("E16").Value = '  -4.04%  '
$ws.Range("D17").Value = '''32.28'
$ws.Range("E17").Value = '  -3.54%  '
$ws.Range("D18").Value = '3.136.48'
$ws.Range("E18").Value = '  -3.95%  '
$ws.Range("E19").Value = '  +4.60%  '
$ws.Range("D20").Value = '''0.0000226'
$ws.Range("E20").Value = '  +19.22%  '
$ws.Range("E21").Value = '  -2.39%  '
$ws.Range("D22").Value = '''426.18'
$ws.Range("E22").Value = '  -2.65%  '
$ws.Range("D23").Value = '''8.38'
$ws.Range("E23").Value = '  -3.36%  '
$ws.Range("E24").Value = '  -4.53%  '
$ws.Range("D25").Value = '''5.44'
$ws.Range("E25").Value = '  +3.88%  '
$ws.Range("D26").Value = '''81.92'
$ws.Range("E26").Value = '  +7.62%  '
$ws.Range("D27").Value = '''11.54'
$ws.Range("E27").Value = '  -3.14%  '
$ws.Range("D28").Value = '3.297.98'
$ws.Range("E28").Value = '  -4.36%  '
$ws.Range("E29").Value = '  +0.02%  '
$ws.Range("D30").Value = '''1.00'
$ws.Range("E30").Value = '  -0.15%  '
$ws.Range("E31").Value = '  -7.34%  '
$ws.Range("E32").Value = '  -4.23%  '
$ws.Range("D33").Value = '''8.20'
$ws.Range("E33").Value = '  -4.24%  '
$ws.Range("D34").Value = '''506.30'
$ws.Range("E34").Value = '  -6.47%  '
$ws.Range("B35").Value = 'Kaspa'
$ws.Range("C35").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D35").Value = '''0.145'
$ws.Range("E35").Value = '  +15.03%  '
$ws.Range("B36").Value = 'RenderToken'
$ws.Range("C36").Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range("D36").Value = '''7.00'
$ws.Range("E36").Value = '  +1.11%  '
$ws.Range("D37").Value = '''1.28'
$ws.Range("E37").Value = '  +2.56%  '
$ws.Range("E38").Value = '  -3.04%  '
$ws.Range("D39").Value = '''22.16'
$ws.Range("E39").Value = '  -0.59%  '
$ws.Range("E40").Value = '  -0.56%  '
$ws.Range("E41").Value = '  -0.07%  '
$ws.Range("E43").Value = '  -3.43%  '
$ws.Range("D44").Value = '''0.364'
$ws.Range("E44").Value = '  -5.01%  '
$ws.Range("D45").Value = '''146.08'
$ws.Range("E45").Value = '  -1.82%  '
$ws.Range("E46").Value = '  +4.90%  '
$ws.Range("D47").Value = '''43.69'
$ws.Range("E47").Value = '  -2.65%  '
$ws.Range("D48").Value = '''0.0674'
$ws.Range("E48").Value = '  +14.39%  '
$ws.Range("D49").Value = '''163.67'
$ws.Range("E49").Value = '  -6.73%  '
$ws.Range("D50").Value = '''0.723'
$ws.Range("E50").Value = '  +1.56%  '
$ws.Range("D51").Value = '''24.15'
$ws.Range("E51").Value = '  -1.48%  '
